$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain-number-looking strings (e.g. "1.004") that must
# stay as TEXT (matching the source inline-string cells), not be auto-converted
# to numbers by Excel. Force text format on those specific cells before writing.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "25.982.44"
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
$ws.Range("D3").Value = "1.639.42"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").Value = "215.03"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("D6").Value = "0.5129"
$ws.Range("E6").Value = "  +1.74%  "

# Row 7
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "0.2584"
$ws.Range("E8").Value = "  +0.44%  "

# Row 9
$ws.Range("D9").Value = "0.06368"
$ws.Range("E9").Value = "  -0.38%  "

# Row 10
$ws.Range("D10").Value = "19.79"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.659.89"
$ws.Range("E12").Value = "  +0.72%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.288"
$ws.Range("E13").Value = "  +0.06%  "

# Row 14
$ws.Range("D14").Value = "0.5469"
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("D15").Value = "64.57"
$ws.Range("E15").Value = "  -0.77%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7738"
$ws.Range("E16").Value = "  -1.73%  "

# Row 17
$ws.Range("D17").Value = "26.009.76"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19
$ws.Range("D19").Value = "197.91"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").Value = "4.438"
$ws.Range("E20").Value = "  +1.16%  "

# Row 21
$ws.Range("D21").Value = "9.974"
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("E22").Value = "  +1.08%  "

# Row 23
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "1.895"
$ws.Range("E24").Value = "  +1.46%  "

# Row 25
$ws.Range("D25").Value = "141.85"
$ws.Range("E25").Value = "  +1.15%  "

# Row 26
$ws.Range("D26").Value = "0.1228"
$ws.Range("E26").Value = "  +7.74%  "

# Row 27
$ws.Range("D27").Value = "6.863"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("D28").Value = "15.67"
$ws.Range("E28").Value = "  -0.36%  "

# Row 29
$ws.Range("D29").Value = "1.239"
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("D30").Value = "0.04886"
$ws.Range("E30").Value = "  -2.34%  "

# Row 31
$ws.Range("D31").Value = "3.280"
$ws.Range("E31").Value = "  +0.61%  "

# Row 32
$ws.Range("E32").Value = "  +0.47%  "

# Row 33
$ws.Range("D33").Value = "1.537"
$ws.Range("E33").Value = "  +0.39%  "

# Row 34
$ws.Range("D34").Value = "2.378"
$ws.Range("E34").Value = "  +0.72%  "

# Row 35
$ws.Range("D35").Value = "0.9134"
$ws.Range("E35").Value = "  +2.14%  "

# Row 36
$ws.Range("D36").Value = "2.591"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
$ws.Range("D37").Value = "0.5534"
$ws.Range("E37").Value = "  +0.20%  "

# Row 38
$ws.Range("D38").Value = "1.114.82"
$ws.Range("E38").Value = "  -2.09%  "

# Row 39
$ws.Range("D39").Value = "0.01568"
$ws.Range("E39").Value = "  +0.88%  "

# Row 40
$ws.Range("D40").Value = "1.004"
$ws.Range("E40").Value = "  -0.01%  "

# Row 41
$ws.Range("D41").Value = "2.540"
$ws.Range("E41").Value = "  -0.65%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.524"
$ws.Range("E42").Value = "  -3.11%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8105"
$ws.Range("E43").Value = "  -0.49%  "

# Row 44
$ws.Range("D44").Value = "99.34"
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("D45").Value = "0.0₈122"
$ws.Range("E45").Value = "  -0.20%  "

# Row 46
$ws.Range("D46").Value = "1.779.91"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47
$ws.Range("D47").Value = "0.4537"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("E48").Value = "  +0.17%  "

# Row 49
$ws.Range("D49").Value = "55.10"
$ws.Range("E49").Value = "  -0.09%  "

# Row 50
$ws.Range("D50").Value = "0.05282"
$ws.Range("E50").Value = "  +4.00%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.09565"
$ws.Range("E51").Value = "  +0.21%  "
